$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.351.67"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "2.284.23"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("D4").Value = '''0.995'
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = '''495.27'
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = '''127.68'
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("D7").Value = '''0.995'
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").Value = '''0.529'
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").Value = "2.282.30"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("D10").Value = '''0.0949'
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "2.660.02"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = '''21.81'
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").Value = "54.185.69"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "2.258.92"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = '''10.05'
$ws.Range("E19").Value = "  +5.41%  "
$ws.Range("D20").Value = '''4.08'
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").Value = '''301.04'
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = '''6.44'
$ws.Range("E22").Value = "  +5.64%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Value = '''62.73'
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").Value = "2.356.56"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").Value = '''169.31'
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "0.0₃0687"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("E40").Value = "  +2.90%  "
$ws.Range("D41").Value = '''3.69'
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = '''0.374'
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("E46").Value = "  +6.34%  "
$ws.Range("D47").Value = '''127.48'
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = '''0.543'
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").Value = '''237.80'
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("D51").Value = '''0.0482'
$ws.Range("E51").Value = "  +2.87%  "
